$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combinations")
$ws.Activate()

# Clear the last data row (row 31): values removed but the row's formatting
# (fill/border styles) stays in place. K31 had no special style, so clearing
# it drops the cell entirely, same as A31:J31 before it had data typed in.
$ws.Range("A31:J31").ClearContents()
$ws.Range("K31").ClearContents()

# Scroll the frozen view down so row 13 is the first visible row under the
# frozen header, then move the selection to B33.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B33").Select()
